# Add data for 2021-12-28: bump the "through December 19" reporting date to
# "through December 20" (sheet name + header label) and update the counts
# that changed as a result of the newly-added day of data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename the sheet / update the report-date label ---------------------
$ws.Name = "Through 2021-12-20"
$ws.Range("B1").Value = "December 2021 (through December 20)"

# --- Update existing counts that changed -----------------------------
$ws.Range("B6").Value = 2
$ws.Range("Z6").Value = 2
$ws.Range("BJ6").Value = 6

$ws.Range("B7").Value = 6
$ws.Range("Z7").Value = 5

$ws.Range("N9").Value = 2

$ws.Range("AX10").Value = 3

$ws.Range("N11").Value = 4
$ws.Range("AL11").Value = 4
$ws.Range("AX11").Value = 4
$ws.Range("BJ11").Value = 5
$ws.Range("BV11").Value = 3

$ws.Range("N15").Value = 6

$ws.Range("AX23").Value = 2

$ws.Range("BV37").Value = 2

$ws.Range("B49").Value = 2

$ws.Range("B54").Value = 3

$ws.Range("B82").Value = 3

# --- New cells (previously empty) ----------------------------------------
$ws.Range("AL4").Value = 1
$ws.Range("BJ16").Value = 1
$ws.Range("BV48").Value = 1
$ws.Range("BJ74").Value = 1
